$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 6 new blank rows at row 198 (this pushes the old rows 198-201 down
# to 204-207, matching the target layout).
# ---------------------------------------------------------------------------
$ws.Rows("198:203").Insert()

# Copy the formatting (styles) of an existing "header + 5 data rows" block
# (rows 165-170: a merged section header followed by 5 REGULAR FOLDER rows
# with the very same B/C/D/E column styles we need) onto the newly inserted
# rows, so the new rows pick up identical cell styles to the target.
$ws.Range("B165:E170").Copy()
$ws.Range("B198:E203").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Re-merge the new section header row (C198:E198), since a format-only paste
# does not restore merged-cell state.
$ws.Range("C198:E198").Merge()

# Match the row height used throughout the rest of the table (15.75pt).
$ws.Range("B198:B203").RowHeight = 15.75

# ---------------------------------------------------------------------------
# Fill in the new section header (row 198).
# ---------------------------------------------------------------------------
$ws.Range("C198").Value = "COMPUTING VOICE LINES (AFTER THE VOICE COMMAND SPOKEN BY THE USER)"

# ---------------------------------------------------------------------------
# Fill in the 5 new data rows (199-203). The "E" column text values are
# entered in this specific order (203, 202, 199, 200, 201) so that the new
# shared-string entries get created in the same order as the source
# workbook (Cyphering, Assessing, Processing, Computing, Evaluating).
# ---------------------------------------------------------------------------
$ws.Range("E203").Value = "Cyphering"
$ws.Range("E202").Value = "Assessing"
$ws.Range("E199").Value = "Processing"
$ws.Range("E200").Value = "Computing"
$ws.Range("E201").Value = "Evaluating"

$ws.Range("B199").Value = "REGULAR FOLDER"
$ws.Range("C199").Value = 20
$ws.Range("D199").Value = "001.wav"

$ws.Range("B200").Value = "REGULAR FOLDER"
$ws.Range("C200").Value = 20
$ws.Range("D200").Value = "002.wav"

$ws.Range("B201").Value = "REGULAR FOLDER"
$ws.Range("C201").Value = 20
$ws.Range("D201").Value = "003.wav"

$ws.Range("B202").Value = "REGULAR FOLDER"
$ws.Range("C202").Value = 20
$ws.Range("D202").Value = "004.wav"

$ws.Range("B203").Value = "REGULAR FOLDER"
$ws.Range("C203").Value = 20
$ws.Range("D203").Value = "005.wav"

# ---------------------------------------------------------------------------
# Match the final cell selection / scroll position recorded in the workbook.
# ---------------------------------------------------------------------------
$null = $ws.Range("E214").Select()
